# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B, C, D, E, G). F (Win) is unchanged.
$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    3 = @(0.0006408296065709695, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.288337117506542)
    4 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    5 = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 27.82738278199502)
    6 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}
